$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.664.35'
$ws.Range("E2").Value = '  +3.27%  '
$ws.Range("D3").Value = '3.072.21'
$ws.Range("E3").Value = '  +3.41%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '518.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.45'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.98%  '
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("E9").Value = '  +1.55%  '
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.376'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.31%  '
$ws.Range("D12").Value = '3.601.60'
$ws.Range("E12").Value = '  +3.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.130'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.76%  '
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '57.723.67'
$ws.Range("E16").Value = '  +3.49%  '
$ws.Range("D17").Value = '3.068.92'
$ws.Range("E17").Value = '  +3.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  +1.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.30%  '
$ws.Range("E25").Value = '  +4.96%  '
$ws.Range("E26").Value = '  +1.20%  '
$ws.Range("D27").Value = '0.0₃0901'
$ws.Range("E27").Value = '  -1.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +3.37%  '
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '154.73'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.95'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.22%  '
$ws.Range("E38").Value = '  +2.08%  '
$ws.Range("D39").Value = '3.114.49'
$ws.Range("E39").Value = '  +3.81%  '
$ws.Range("E40").Value = '  +4.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.654'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '
$ws.Range("D44").Value = '2.263.89'
$ws.Range("E44").Value = '  +4.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0259'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '20.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.90%  '
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.922'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '261.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.711'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.00%  '
